$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values are plain decimal numbers (e.g. "2.00", "9.50").
# Excel would auto-convert these to numeric cells (dropping trailing zeros / exact
# text form) unless the cell is pre-formatted as Text, so force Text format first
# to preserve them as literal strings, matching the source data feed formatting.
$textCells = @("D4", "D6", "D8", "D12", "D19", "D20", "D21", "D23", "D25", "D26", "D27", "D39", "D40", "D41", "D43", "D46", "D48", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values (coin names, links, prices, 1h volume deltas).
$ws.Range("D2").Value = "62.398.99"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.449.89"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "144.37"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "2.443.53"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "62.182.42"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.447.80"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "10.92"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "7.17"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "330.51"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "66.11"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("D27").Value = "628.42"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.564.42"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "5.35"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "149.91"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "18.38"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "42.47"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "143.97"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").Value = "0.0526"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "0.601"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "19.63"
$ws.Range("E50").Value = "  -7.03%  "
$ws.Range("D51").Value = "0.0₆0239"
$ws.Range("E51").Value = "  +8.89%  "
